# Weekly update: insert 4 new daily-price rows for Jengibre (Vega Central
# Mapocho de Santiago) at the top of the data block (rows 28-31), pushing the
# existing historical rows (old 28-47) down to rows 32-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 28 (shifts old rows 28-47 -> 32-51)
$ws.Range("A28:A31").EntireRow.Insert()

# Common (constant) values shared by every data row in this sheet
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100114007
$categoria = "Jengibre"
$variedad  = "Sin especificar"
$origen    = "Perú"
$clasif    = "Hortaliza"

function Set-JengibreRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PMin,
        [double]$PMax,
        [double]$PProm,
        [string]$Unidad,
        [double]$PrecioKg,
        [double]$KgUnidades
    )

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $catId
    $ws.Cells.Item($Row, 7).Value  = $categoria
    $ws.Cells.Item($Row, 8).Value  = $variedad
    $ws.Cells.Item($Row, 9).Value  = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PMin
    $ws.Cells.Item($Row, 12).Value = $PMax
    $ws.Cells.Item($Row, 13).Value = $PProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasif
}

# New row 28: 2021-08-27, Primera
Set-JengibreRow 28 44435 "Primera" 880 13000 14000 13500 "$/caja 13 kilos" 1038 13

# New row 29: 2021-08-27, Segunda
Set-JengibreRow 29 44435 "Segunda" 340 11000 12000 11500 "$/caja 13 kilos" 885 13

# New row 30: 2021-08-23, Primera
Set-JengibreRow 30 44431 "Primera" 880 13000 14000 13500 "$/caja 13 kilos" 1038 13

# New row 31: 2021-08-23, Segunda
Set-JengibreRow 31 44431 "Segunda" 340 11000 12000 11500 "$/caja 13 kilos" 885 13
